$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header for column B from "Ename" to "Name"
$ws.Range("B3").Value = "Name"

# Column F (BgImage) rows 4-7: change numeric level values to the text "t5"
$ws.Range("F4").Value = "t5"
$ws.Range("F5").Value = "t5"
$ws.Range("F6").Value = "t5"
$ws.Range("F7").Value = "t5"

# Update the active selection to match the edited region
$ws.Range("F5:F7").Select()
